# Convert adoc tables to html tables
# Prefix each table row with a leading "|" (adoc -> html table syntax),
# and fix a couple of small data typos in the "many-rows" example table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# table1 (Rule Description for table1 anchors-in-cells example)
# Leading apostrophe forces this to be entered as literal text rather
# than being parsed as a formula, since the text starts with "==".
$ws.Range("C27").Value = "'===`n|WITH anchor`n|WITHOUT anchor`n==="

# table2 (Rule Description for table2 example)
$ws.Range("C28").Value = "|Header 1|Header 2`n===`n|Cell in column 1, row 1|Cell in column 2, row 1`n|Cell in column 1, row 2|Cell in column 2, row 2`n==="

# table5 (Rule Description for the many-rows example; also corrects
# Color8 -> Color7 typo and adds the missing Name8|Color8 row while
# dropping the trailing Name10|Color10 row)
$ws.Range("C31").Value = "|Name|Color`n===`n|Roses|Red`n|Violets|Blue`n|Name1|Color1`n|Name2|Color2`n|Name3|Color3`n|Name4|Color4`n|Name5|Color5`n|Name6|Color6`n|Name7|Color7`n|Name8|Color8`n|Name9|Color9`n..."
